$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet2 ("suite-demo1"): replace the old "testRetailHomepage" test case with
# a new, smaller "testEmailServices" test case.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("suite-demo1")

# Old layout occupied rows 2-8 (row1 was blank). Shift everything up by one
# so the header row becomes row 1.
$ws2.Rows("1:1").Delete()

# Rows are now: 1=header, 2=open/domain, 3=click/searchBox, 4=type/..,
# 5=click/searchButton, 6=waitForPageToLoad/10, 7=assertTextPresent/..
# Drop the trailing row (7) completely.
$ws2.Rows("7:7").Delete()
# Blank out rows 4 and 5 entirely (no cells left at all in target).
$ws2.Rows("4:5").Clear()
# Row 6 keeps only an (empty, styled) C6 cell - drop its old B6/C6 values.
$ws2.Range("B6").ClearContents()
$ws2.Range("C6").ClearContents()

# New header row.
$ws2.Range("A1").Value = "testEmailServices"
$ws2.Range("B1").Value = "Function"
$ws2.Range("C1").Value = "Target"
$ws2.Range("D1").Value = "Value"

# New test steps.
$ws2.Range("B2").Value = "open"
$ws2.Range("C2").Value = "/email-services/"
$ws2.Range("B3").Value = "matt"
$ws2.Range("C3").Value = "email"

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws2.Range("A1:XFD13").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet3 ("data set 1"): append a new "message"/"domain" data row.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("data set 1")
$ws3.Activate()
$ws3.Range("A7").Value = "message"
$ws3.Range("B7").Value = "domain"
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1
$ws3.Range("B8").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet4 ("data set 2"): append a new "message"/"domain" data row.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("data set 2")
$ws4.Activate()
$ws4.Range("A7").Value = "message"
$ws4.Range("B7").Value = "domain"
$ws4.PageSetup.PaperSize = 9
$ws4.PageSetup.Orientation = 1
$ws4.Range("B8").Select() | Out-Null

# Restore "suite-demo1" as the active/selected sheet+tab.
$ws2.Activate()
